$d = $word.ActiveDocument

# The last paragraph currently holds only the "_GoBack" bookmark and is
# center-aligned. Drop the centering (the new paragraphs / run are left
# aligned with bold text instead), then insert two new bold paragraphs
# ("Part I", "Part II") right before it, and finally insert a bold
# "Part III" run into the bookmark paragraph itself, ahead of the bookmark.

$lastPara = $d.Paragraphs.Last
$lastPara.Alignment = 0

$r = $lastPara.Range
$r.InsertParagraphBefore() | Out-Null
$partI = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$partI.Range.InsertAfter("Part I: MMM using global memory only")
$partI.Range.Bold = 1

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphBefore() | Out-Null
$partII = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$partII.Range.InsertAfter("Part II: MMM using shared memory")
$partII.Range.Bold = 1

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertBefore("Part III: MMM using shared memory and loop unrolling")
$lastPara = $d.Paragraphs.Last
$lastPara.Range.Bold = 1
